# Trade #112 closed at 2026-02-16 21:43:05 - momentum DOWN +0.000%
#
# This applies four groups of edits:
#  1. Summary sheet    - refreshed aggregate stats (trade #88 now closed)
#  2. Comparison sheet - refreshed momentum-strategy stats
#  3. momentum sheet   - trade #88 (row 23) flips from OPEN -> CLOSED,
#                        and a brand new trade #112 (row 27) is appended as OPEN
#  4. All Trades sheet - the now-closed trade #88 is appended as a new row (89)
#
# NOTE: several source values look numeric/date-like ("70.5%", "11.39",
# "2026-02-16", ...) but must stay literal text, matching every other cell
# in these report sheets (they're pre-formatted strings, not real
# Excel numbers/dates). Excel's auto-detection would otherwise coerce them
# into a percentage/number/date, so NumberFormat is forced to Text ("@")
# immediately before those particular writes.

$wb = $excel.ActiveWorkbook

function Set-TextValue($range, $value) {
    $range.NumberFormat = "@"
    $range.Value = $value
}

# ---------------------------------------------------------------------------
# 1. Summary sheet
# ---------------------------------------------------------------------------
$summary = $wb.Worksheets.Item("Summary")

$summary.Range("C2").Value = 88
Set-TextValue $summary.Range("D2") "70.5%"
Set-TextValue $summary.Range("E2") "+25.8308%"
Set-TextValue $summary.Range("F2") "+0.2935%"

Set-TextValue $summary.Range("D4") "80.0%"
Set-TextValue $summary.Range("E4") "+11.6836%"
Set-TextValue $summary.Range("F4") "+0.4673%"

# ---------------------------------------------------------------------------
# 2. Comparison sheet
# ---------------------------------------------------------------------------
$comparison = $wb.Worksheets.Item("Comparison")

Set-TextValue $comparison.Range("C3") "80.0%"
Set-TextValue $comparison.Range("D3") "11.39"
Set-TextValue $comparison.Range("E3") "+0.6404%"
Set-TextValue $comparison.Range("G3") "1.14"

# ---------------------------------------------------------------------------
# 3. momentum sheet
# ---------------------------------------------------------------------------
$momentum = $wb.Worksheets.Item("momentum")

# Trade #88 (row 23) closes out
$momentum.Range("G23").Value = 68233.05665499999
$momentum.Range("H23").Value = "CLOSED"
$momentum.Range("I23").Value = 0.3623
$momentum.Range("J23").Value = 3.62
$momentum.Range("M23").Value = "time_exit_5min"
$momentum.Range("N23").Value = 5

# New trade #112 (row 27), freshly opened
$momentum.Cells.Item(27, 1).Value = 112
Set-TextValue $momentum.Cells.Item(27, 2) "2026-02-16"
$momentum.Cells.Item(27, 3).Value = "21:43:05"
$momentum.Cells.Item(27, 4).Value = "momentum"
$momentum.Cells.Item(27, 5).Value = "DOWN"
$momentum.Cells.Item(27, 6).Value = 68355.74000000001
$momentum.Cells.Item(27, 8).Value = "OPEN"
$momentum.Cells.Item(27, 9).Value = 0
$momentum.Cells.Item(27, 10).Value = 0
$momentum.Cells.Item(27, 11).Value = 0.9
$momentum.Cells.Item(27, 12).Value = "Downward momentum: -0.175% over 10 samples"
$momentum.Cells.Item(27, 14).Value = 0

# ---------------------------------------------------------------------------
# 4. All Trades sheet - append the now-closed trade #88 as a new row
# ---------------------------------------------------------------------------
$allTrades = $wb.Worksheets.Item("All Trades")

$allTrades.Cells.Item(89, 1).Value = 88
Set-TextValue $allTrades.Cells.Item(89, 2) "2026-02-16"
$allTrades.Cells.Item(89, 3).Value = "21:38:01"
$allTrades.Cells.Item(89, 4).Value = "momentum"
$allTrades.Cells.Item(89, 5).Value = "DOWN"
$allTrades.Cells.Item(89, 6).Value = 68481.19500000001
$allTrades.Cells.Item(89, 7).Value = 68233.05665499999
$allTrades.Cells.Item(89, 8).Value = "CLOSED"
$allTrades.Cells.Item(89, 9).Value = 0.3623
$allTrades.Cells.Item(89, 10).Value = 3.62
$allTrades.Cells.Item(89, 11).Value = 0.9
$allTrades.Cells.Item(89, 12).Value = "Downward momentum: -0.132% over 10 samples"
$allTrades.Cells.Item(89, 13).Value = "time_exit_5min"
$allTrades.Cells.Item(89, 14).Value = 5
